$d = $word.ActiveDocument

$found = $d.Content.Find.Execute("Matching product item on each location", $true, $false, $false, $false, $false, $true, 1, $false, "Matching product on each location", 2)
Write-Host "Find1: $found"
